$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status moved from "Ready for handoff" to "In Translation" for both locales
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# The Status columns narrow to reflect the shorter text (was autofit to
# "Ready for handoff", now autofit to "In Translation")
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
